$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "26 марта"
$ws.Range("B31").Value = "Декомпозиция кода (написание NetworkReceiveManager и вынесение структуры Measurements)"
$ws.Range("B31").HorizontalAlignment = -4152

$ws.Range("A32").Value = "27 марта"
$ws.Range("B32").Value = "Написание генератора для выдачи 4к изображений из видео и адаптация и мелкие фиксы"
$ws.Range("B32").HorizontalAlignment = -4152

$ws.Range("B33").Select()
